# chore(auto): pipeline update & logic validation
# Refresh the war-participation report: the data pipeline re-pulled the
# clan roster, which re-sorted the existing members and added one new
# member ("StelaAby"). Rewrite the full data block (rows 2..65) in the
# new order/content, then leave the sheet's used range at A1:H65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("JOSA","Razoável","Guerra Atual",12,10,16,11,13),
  @("nova","Verificar","Guerra Atual",0,0,0,0,0),
  @("Yzzy the best ¿","Verificar","Guerra Atual",0,0,0,0,0),
  @("Gustavo","Verificar","Guerra Atual",0,0,0,0,0),
  @("sj pro","Verificar","Guerra Atual",0,8,10,0,0),
  @("alm ícaro⁷⁷⁷","Verificar","Guerra Atual",0,7,0,0,0),
  @("- ROJÎ","Verificar","Guerra Atual",0,0,0,0,0),
  @("Felipe","Verificar","Guerra Atual",0,0,0,0,0),
  @("Zenitsu愛","Verificar","Guerra Atual",3,5,12,12,12),
  @("danilo","Ok","Guerra Atual",16,16,16,5,0),
  @("PedrinhoR14","Razoável","Guerra Atual",14,16,16,16,16),
  @("5C4RF4C3","Verificar","Guerra Atual",0,0,0,0,0),
  @("luba","Verificar","Guerra Atual",0,12,8,16,16),
  @("ed","Ok","Guerra Atual",16,16,16,16,16),
  @("caioba_","Verificar","Guerra Atual",0,1,0,0,0),
  @("BRS⚔️ASHURA","Verificar","Guerra Atual",9,16,16,16,12),
  @("john","Verificar","Guerra Atual",0,15,9,11,12),
  @("Theus Carvalho","Razoável","Guerra Atual",12,11,16,14,16),
  @("mini.peka","Verificar","Guerra Atual",4,0,0,0,0),
  @("kauansin777","Verificar","Guerra Atual",0,12,12,12,9),
  @("Luiz Fernando™","Verificar","Guerra Atual",8,16,12,13,10),
  @("dogmal","Verificar","Guerra Atual",0,0,0,0,0),
  @("StelaAby","Verificar","Guerra Atual",0,0,0,0,0),
  @("^_^^_^","Razoável","Guerra Atual",12,16,12,12,16),
  @("WILLIAN","Verificar","Guerra Atual",0,0,0,0,0),
  @("Chetto","Razoável","Guerra Atual",12,13,11,16,15),
  @("Alvaro","Razoável","Guerra Atual",13,14,15,14,14),
  @("Nico","Razoável","Guerra Atual",14,14,9,14,16),
  @("Pedro PH","Ok","Guerra Atual",16,16,8,0,0),
  @("polaris","Ok","Guerra Atual",16,16,16,16,16),
  @("GabiMalvadeza","Verificar","Guerra Atual",0,0,0,0,0),
  @("tavin","Razoável","Guerra Atual",12,8,14,12,16),
  @("Grimmer 狼","Ok","Guerra Atual",16,4,0,0,0),
  @("Rodolfos","Ok","Guerra Atual",16,16,16,15,12),
  @("super","Ok","Guerra Atual",16,16,16,16,16),
  @("domador de but","Ok","Guerra Atual",16,16,8,12,12),
  @("gabiggoughost","Ok","Guerra Atual",16,14,16,16,16),
  @("Daniele❤","Ok","Guerra Atual",16,16,16,16,16),
  @("Teixeirazzqw","Verificar","Guerra Atual",10,0,0,0,0),
  @("Luciano","Ok","Guerra Atual",16,16,14,16,16),
  @("filho de duque","Verificar","Guerra Atual",0,0,0,0,0),
  @("Sotto ツ","Razoável","Guerra Atual",12,14,12,13,16),
  @("OneDePrata","Razoável","Guerra Atual",15,8,16,13,16),
  @("mathi❄️","Razoável","Guerra Atual",14,10,16,16,16),
  @("Gustavo Clash","Ok","Guerra Atual",16,16,16,16,16),
  @("Dockz","Ok","Guerra Atual",16,14,16,16,15),
  @("RaiNascimento","Verificar","Guerra Atual",0,0,0,0,0),
  @("luck","Ok","Guerra Atual",16,16,12,16,13),
  @("Bruno","Ok","Guerra Atual",16,16,0,0,16),
  @("andrebts","Razoável","Guerra Atual",15,16,16,16,16),
  @("⭐O SENTINELA ⭐","Ok","Guerra Atual",16,16,16,16,16),
  @("O GUARDIÃO","Ok","Guerra Atual",16,16,16,16,16),
  @("DGJ-DAVI","Ok","Guerra Atual",16,16,12,15,16),
  @("51 é pinga","Ok","Guerra Atual",16,16,16,16,16),
  @("EDDIE","Ok","Guerra Atual",16,16,16,16,16),
  @("WvCly","Verificar","Guerra Atual",8,14,7,10,6),
  @("isp","Ok","Guerra Atual",16,16,16,16,16),
  @("Mila","Verificar","Guerra Atual",8,16,16,16,16),
  @("gabriel 3`$","Ok","Guerra Atual",16,16,16,16,4),
  @("nivelador","Razoável","Guerra Atual",12,10,15,12,0),
  @("Asten Acady","Ok","Guerra Atual",16,16,16,16,16),
  @("joão3:16","Ok","Guerra Atual",16,16,16,16,16),
  @("ZackThunder","Ok","Guerra Atual",16,8,16,16,15),
  @("RobaFrag","Ok","Guerra Atual",16,16,16,16,16)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
